# Refresh the cryptocurrency price ("Price", column D) and 1h volume-change
# ("Volume(1h)", column E) figures on the sheet with the latest scraped
# values. Both columns hold plain text (e.g. prices can contain more than
# one "." like "56.249.09", and the volume column keeps its padding spaces
# around the percentage, e.g. "  -3.41%  "), so every write that could be
# misread by Excel as a genuine number first forces the cell to Text format
# and restores the original (Normal) style afterwards, leaving appearance
# untouched while keeping the value stored as a string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '56.249.09'
$ws.Range("E2").Value = '  -3.41%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.376.84'
$ws.Range("E3").Value = '  -3.32%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.06%  '

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '500.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.17%  '

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.95'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.69%  '

# Row 7: USDC
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.40%  '

# Row 8: XRP
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.553'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.53%  '

# Row 9: LidoStakedEther
$ws.Range("D9").Value = '2.398.89'
$ws.Range("E9").Value = '  -2.41%  '

# Row 10: Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0955'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.47%  '

# Row 11: TRON
$ws.Range("E11").Value = '  -1.02%  '

# Row 12: Toncoin
$ws.Range("E12").Value = '  -7.17%  '

# Row 13: Cardano
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.316'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.34%  '

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = '2.803.13'
$ws.Range("E14").Value = '  -3.15%  '

# Row 15: WrappedBTC
$ws.Range("D15").Value = '56.124.02'
$ws.Range("E15").Value = '  -3.46%  '

# Row 16: Avalanche
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.43'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.19%  '

# Row 17: ShibaInu
$ws.Range("E17").Value = '  -1.58%  '

# Row 18: WrappedEther
$ws.Range("D18").Value = '2.407.44'
$ws.Range("E18").Value = '  -2.26%  '

# Row 19: Chainlink
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.07'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.48%  '

# Row 20: BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '308.99'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.67%  '

# Row 21: Polkadot
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.01'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.98%  '

# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.26'
$ws.Range("D22").Style = "Normal"

# Row 24: LEO
$ws.Range("E24").Value = '  -4.57%  '

# Row 25: Litecoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.44'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.09%  '

# Row 26: Binance-PegBSC-USD
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.997'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.70%  '

# Row 27: WrappedeETH
$ws.Range("D27").Value = '2.486.85'
$ws.Range("E27").Value = '  -3.85%  '

# Row 28: Polygon
$ws.Range("E28").Value = '  -6.92%  '

# Row 29: Kaspa
$ws.Range("E29").Value = '  -5.41%  '

# Row 30: InternetComputer(DFINITY)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.25'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.60%  '

# Row 31: Monero
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '172.51'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.88%  '

# Row 32: PEPE
$ws.Range("D32").Value = '0.0₃0711'
$ws.Range("E32").Value = '  -3.75%  '

# Row 33: PancakeSwap
$ws.Range("E33").Value = '  -2.68%  '

# Row 34: Aptos
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.09'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.68%  '

# Row 35: USDe
$ws.Range("E35").Value = '  -0.10%  '

# Row 36: Fetch.AI
$ws.Range("E36").Value = '  -6.39%  '

# Row 37: FirstDigitalUSD
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.994'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.54%  '

# Row 38: EthereumClassic
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.83'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.01%  '

# Row 39: ImmutableX
$ws.Range("E39").Value = '  +0.92%  '

# Row 40: NEARProtocol
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.76'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.68%  '

# Row 41: OKB
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '35.82'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.16%  '

# Row 42: Stacks
$ws.Range("E42").Value = '  -3.04%  '

# Row 43: SuiNetwork
$ws.Range("E43").Value = '  -1.98%  '

# Row 44: Aave
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '129.12'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.39%  '

# Row 45: Filecoin
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.33'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.51%  '

# Row 46: RenderToken
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.73'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.64%  '

# Row 47: Bittensor
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '250.69'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.68%  '

# Row 48: Mantle
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.560'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.01%  '

# Row 49: Stellar
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0899'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.15%  '

# Row 50: Hedera
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0483'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.03%  '

# Row 51: EnergySwap
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.84'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.62%  '

